# Elimina antiguos EC y agrega nuevos y modifica Antigua BD
#
# Functional change: the "Periodo Mora" value shown for every worker
# (rows 16-19, column E) moves from period 2508 to period 2509, and
# those four cells pick up center horizontal alignment (matching the
# rest of the data rows/table).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Periodo Mora" text for each worker row (2508 -> 2509) and
# center it horizontally, same as the other data columns in the table.
foreach ($r in 16..19) {
    $cell = $ws.Cells.Item($r, 5)   # column E
    $cell.Value = "2509"
    $cell.HorizontalAlignment = -4108   # xlCenter
}
